$d = $word.ActiveDocument

# Locate the paragraph containing the old single-run text and rewrite it in
# place (this keeps the existing run/paragraph formatting intact).
$rng = $d.Content
$old = "Git –ammend  for editing the commit msg"
$new = "Git –ammend  for editing the commit message"
$found = $rng.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)

if ($found) {
    $start = $rng.Start
    $splitAt = $start + 7  # length of "Git –am"
    $end = $start + $new.Length

    # Re-apply (and then revert) a character formatting property on the
    # first chunk so the run splits into two separate <w:r> elements with
    # matching rPr, mirroring the edit's two-run result.
    $r1 = $d.Range($start, $splitAt)
    $r1.Bold = 1
    $r1.Bold = 0

    $r2 = $d.Range($splitAt, $end)
    Write-Output ("r1=[" + $r1.Text + "] r2=[" + $r2.Text + "]")
}
else {
    Write-Output "Target text not found"
}
